$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Combine the old row 2..5 contents ("Elemental", "Token Creature — Elemental",
# "This creature’s power and toughness are each equal to the number of
# creatures you control.", "*/*") into a single Python-tuple-looking string
# and store it in A2.
$ws.Range("A2").Value = "('Elemental', ['Token Creature — Elemental', 'This creature’s power and toughness are each equal to the number of creatures you control.', '*/*'])"

# The old rows 3-5 are no longer needed now that their data lives in A2.
$ws.Range("A3:A5").EntireRow.Delete()
